# Apply updated "想去人数" (F column) counts scraped at commit 456a3b4.
$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 2612
$ws.Range("F5").Value = 1315
$ws.Range("F7").Value = 3197
$ws.Range("F11").Value = 8132
$ws.Range("F19").Value = 498
$ws.Range("F23").Value = 280
$ws.Range("F40").Value = 2341
$ws.Range("F41").Value = 1217

# Sheet "演出"
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F5").Value = 174

# Sheet "全部类型"
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 2612
$ws.Range("F4").Value = 174
$ws.Range("F6").Value = 1315
$ws.Range("F8").Value = 3197
$ws.Range("F13").Value = 8132
$ws.Range("F20").Value = 498
$ws.Range("F23").Value = 280
$ws.Range("F39").Value = 2341
$ws.Range("F41").Value = 1217
